$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 19
$values = @(2026, 2, 22, 5102, 4924, 4843, 4768, 4742, 4782, 4919, 5216, 5593, 5872, 5921, 5857, 5769, 5708, 5554, 5502, 5734, 6126, 6611, 6614, 6452, 6207, 6014, 5738)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($row, $i + 1).Value = $values[$i]
}
